$d = $word.ActiveDocument

# --- Locate the paragraph that ends with "...front-end and back-end development!" ---
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*front-end and back-end development!*") {
        $targetPara = $p
    }
}

$paraRange = $targetPara.Range
# Position right before the paragraph mark (i.e. end of the visible text).
$textEnd = $paraRange.End - 1

$newText = " My goal is to eventually land a job as a junior developer on a great time. I" + [char]0x2019 + "m committed to a lifetime of learning in this field to become the best I can be."

# Insert the new sentence at the end of the paragraph. Because the insertion
# point is collapsed right after the existing text, it inherits that run's
# character formatting (Calibri Light, 12pt, incl. complex-script font).
$insertRange = $d.Range($textEnd, $textEnd)
$insertRange.InsertAfter($newText)

$newRunStart = $textEnd
$newRunEnd = $textEnd + $newText.Length

# Force the freshly typed text into its own <w:r> (the diff adds a distinct
# run rather than extending the previous one) by toggling a boolean character
# property on and back off right over that span.
$newRunRange = $d.Range($newRunStart, $newRunEnd)
$newRunRange.Bold = 1
$newRunRange.Bold = 0

# --- Move the _GoBack bookmark from the end of the document to right after
# the sentence we just added (still inside this paragraph, before its mark).
#
# Quirk workaround: adding a zero-length bookmark exactly at a paragraph-mark
# boundary confuses this host, so we temporarily insert a one-character
# marker past the target spot (pushing that boundary away), add the bookmark,
# then delete the marker again. The bookmark stays put.
$marker = $d.Range($newRunEnd, $newRunEnd)
$marker.InsertAfter("Z")

$bookmarkRange = $d.Range($newRunEnd, $newRunEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$markerRange = $d.Range($newRunEnd, $newRunEnd + 1)
$markerRange.Delete()
